$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update "Valor Mora" total
$ws.Range("E11").Value = 409066

# Replace the worker/period table (rows 16-23) with the updated dataset
$data = @(
    @("CC", "1069473966", "JORGE LUIS RUIZ SOTO", "2412", 52000, 1300000),
    @("CC", "1051893091", "DAILYN PEREZ BLANCO", "2412", 52000, 1300000),
    @("CC", "1069473966", "JORGE LUIS RUIZ SOTO", "2501", 52000, 1300000),
    @("CC", "1069473966", "JORGE LUIS RUIZ SOTO", "2502", 52000, 1300000),
    @("CC", "1069473966", "JORGE LUIS RUIZ SOTO", "2503", 52000, 1300000),
    @("CC", "1069473966", "JORGE LUIS RUIZ SOTO", "2504", 52000, 1300000),
    @("CC", "1069473966", "JORGE LUIS RUIZ SOTO", "2505", 52000, 1300000),
    @("CC", "1069473966", "JORGE LUIS RUIZ SOTO", "2506", 45066, 1300000)
)

$r = 16
foreach ($row in $data) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}

# Remove the two now-unused table rows (24 and 25); this shifts the
# footer block (rows 30-31) up to rows 28-29, matching the new layout.
$ws.Range("A24:A25").EntireRow.Delete()
